$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '31.393.85'
$ws.Range("E2").Value = '  +3.48%  '

$ws.Range("D3").Value = '1.988.45'
$ws.Range("E3").Value = '  +6.28%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.28%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.8041'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +71.40%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '252.81'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.78%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.003'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.25%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3414'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +18.78%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.54'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +15.90%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06970'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.04%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8405'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +16.54%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08120'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.45%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '101.35'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.01%  '

$ws.Range("D14").Value = '1.989.44'
$ws.Range("E14").Value = '  +6.28%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.454'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.38%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '272.80'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.20%  '

$ws.Range("D17").Value = '31.389.79'
$ws.Range("E17").Value = '  +3.49%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.93'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +7.29%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007946'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.64%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.725'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +9.18%  '

$ws.Range("D21").Value = '2.250.25'
$ws.Range("E21").Value = '  +6.37%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.15%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.004'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.39%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.965'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +11.61%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.660'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.81%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1524'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +59.11%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '165.37'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.74'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.76%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.178'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +16.21%  '

$ws.Range("E30").Value = '  +6.73%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.358'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.78%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.581'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +8.98%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.331'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.69%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05181'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.85%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.213'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.41%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7501'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +8.89%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.811'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.85%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01998'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.40%  '

$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.928'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.21%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.631'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.69%  '

$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '78.83'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.27%  '

$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4651'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +10.17%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.060'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.08%  '

$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '105.87'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.20%  '

$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8533'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.48%  '

$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.003'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.37%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.976'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.20%  '

$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.482'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +8.11%  '

$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.54'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.88%  '

$ws.Range("B50").Value = 'Decentraland'
$ws.Range("C50").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4275'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +9.16%  '

$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1179'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +10.39%  '
